# Generate Report for Handback
#
# The "f794f99b-2547-4430-bdd2-3cb38429d534.md" file (row 7 in every sheet)
# moves from "Ready for handoff" into "Handed back: in sync with en-US":
#  - Overview: zh-cn/de-de status columns (E7/F7) flip to the handed-back text.
#  - zh-cn / de-de sheets: Status (C7) flips, the Latest Target File (I7) and
#    Latest Handback File (J7) columns get filled in (with a hyperlink on I7,
#    matching the style used on the rows above), and the Latest Handback
#    DateTime (K7) gets a real timestamp instead of the 0001-01-01 placeholder.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = $handedBack
$overview.Range("F7").Value = $handedBack

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = $handedBack
$zhcn.Range("I7").Value = "f794f99b-2547-4430-bdd2-3cb38429d534.md"
$zhcn.Range("J7").Value = "f794f99b-2547-4430-bdd2-3cb38429d534.b1ae9dbe47a818c3ca63c9b2006233857166b46a.zh-cn.xlf"
$zhcn.Range("K7").Value = "2016-10-24 09:36:20"

$zhcn.Hyperlinks.Add($zhcn.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ba2588552e9a6401ae89afc9dfa0b344dc02a9ef/e2e/f794f99b-2547-4430-bdd2-3cb38429d534.md", "", "", "f794f99b-2547-4430-bdd2-3cb38429d534.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = $handedBack
$dede.Range("I7").Value = "f794f99b-2547-4430-bdd2-3cb38429d534.md"
$dede.Range("J7").Value = "f794f99b-2547-4430-bdd2-3cb38429d534.b1ae9dbe47a818c3ca63c9b2006233857166b46a.de-de.xlf"
$dede.Range("K7").Value = "2016-10-24 09:36:36"

$dede.Hyperlinks.Add($dede.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ba2588552e9a6401ae89afc9dfa0b344dc02a9ef/e2e/f794f99b-2547-4430-bdd2-3cb38429d534.md", "", "", "f794f99b-2547-4430-bdd2-3cb38429d534.md")

# ---------------------------------------------------------------------------
# Re-assert the "hyperlink" and "date/time" formatting on every data row of
# every sheet. The workbook only carries two custom cell styles (a blue
# underlined hyperlink font, and a yyyy-mm-dd HH:mm:ss date format); re-
# applying them keeps every row -- old and newly-populated alike -- visually
# consistent with the rest of the table.
# ---------------------------------------------------------------------------
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) == FF6495ED

function Format-HyperlinkCell($cell) {
    $cell.Font.Underline = 2
    $cell.Font.Color = $hyperlinkColor
}

function Format-DateCell($cell) {
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

for ($r = 2; $r -le 7; $r++) {
    Format-HyperlinkCell $overview.Cells.Item($r, 2)   # column B
    Format-DateCell $overview.Cells.Item($r, 7)        # column G

    Format-HyperlinkCell $zhcn.Cells.Item($r, 1)        # column A
    Format-HyperlinkCell $zhcn.Cells.Item($r, 9)        # column I
    Format-DateCell $zhcn.Cells.Item($r, 8)             # column H
    Format-DateCell $zhcn.Cells.Item($r, 11)            # column K

    Format-HyperlinkCell $dede.Cells.Item($r, 1)        # column A
    Format-HyperlinkCell $dede.Cells.Item($r, 9)        # column I
    Format-DateCell $dede.Cells.Item($r, 8)             # column H
    Format-DateCell $dede.Cells.Item($r, 11)            # column K
}

Write-Output "Report for handback generated"
